$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

# Paragraph 9 currently holds "@Tag" - split it so a new, empty paragraph
# follows it, then fill that new paragraph with the run-by-run text of the
# "API Assertivas - Hamcrest and AssertJ" line.
$tagPara = $tr.Paragraphs(9, 1)
[void]$tagPara.InsertAfter([char]13)

$cur = $tr.Paragraphs(10, 1)
$cur = $cur.InsertAfter("API Assertivas - ")
$cur = $tr.Paragraphs(10, 1).InsertAfter("Hamcrest")
$cur = $tr.Paragraphs(10, 1).InsertAfter(" ")
$cur = $tr.Paragraphs(10, 1).InsertAfter("and")
$cur = $tr.Paragraphs(10, 1).InsertAfter(" ")
$cur = $tr.Paragraphs(10, 1).InsertAfter("AssertJ")
